$wb = $excel.ActiveWorkbook

$newA8 = "S1" + [char]10 + "(07:00-09:00)"
$newB8 = "Võ Văn F"
$newD8 = "Lớp: CL05" + [char]10 + "Môn: Tiếng Anh chuyên ngành" + [char]10 + "Phòng: R101" + [char]10 + "(Lý thuyết)"

$newA9 = "C2" + [char]10 + "(15:00-17:00)"
$newB9 = "Hoàng Thị E"
$newG9 = "Lớp: CL05" + [char]10 + "Môn: Kỹ năng mềm" + [char]10 + "Phòng: R101" + [char]10 + "(Lý thuyết)"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # ----- Row 8 -----
    # A8: update time slot label (style 6 unchanged)
    $ws.Range("A8").Value = $newA8

    # B8: update teacher name (style 7 unchanged)
    $ws.Range("B8").Value = $newB8

    # D8: becomes the highlighted "occupied" cell (style 8), carrying the class info
    #     that used to live in E8. Copy E8's format (style 8) into D8 first, then set
    #     the new text.
    $ws.Range("E8").Copy()
    $ws.Range("D8").PasteSpecial(-4122)
    $ws.Range("D8").Value = $newD8

    # E8: no longer highlighted, becomes a plain bordered cell (style 7), with no content.
    $ws.Range("E8").ClearContents()
    $ws.Range("C8").Copy()
    $ws.Range("E8").PasteSpecial(-4122)

    # ----- Row 9 -----
    # A9: was an empty numeric cell (style 9); becomes a time slot label cell like A8
    #     (style 6), and the row grows to the same 60pt height as row 8.
    $ws.Range("A8").Copy()
    $ws.Range("A9").PasteSpecial(-4122)
    $ws.Range("A9").Value = $newA9
    $ws.Rows.Item(9).RowHeight = 60

    # B9: update teacher name (style 7 unchanged)
    $ws.Range("B9").Value = $newB9

    # G9: becomes the highlighted "occupied" cell (style 8), carrying the class info
    #     that used to live in F9. Copy F9's format (style 8) into G9 first, then set
    #     the new text.
    $ws.Range("F9").Copy()
    $ws.Range("G9").PasteSpecial(-4122)
    $ws.Range("G9").Value = $newG9

    # F9: no longer highlighted, becomes a plain bordered cell (style 7), with no content.
    $ws.Range("F9").ClearContents()
    $ws.Range("C9").Copy()
    $ws.Range("F9").PasteSpecial(-4122)
}
